$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93: the Publicado_em_BRT (B) cell was an empty placeholder that gets
# dropped entirely (no value was ever published for this item).
$ws.Range("B93").ClearContents()

# Row 94: new news item appended by the bot run.
$ws.Range("A94").Value = "05/01/2026 14:21:04"
$ws.Range("B94").Value = "05/01 14:14"
$ws.Range("C94").Value = "g1 > Política"
$ws.Range("D94").Value = "Sarney critica invasão na Venezuela e elogia posição do Brasil; VEJA VÍDEO"
$ws.Range("E94").Value = "https://g1.globo.com/politica/blog/gerson-camarotti/post/2026/01/05/sarney-critica-invasao-na-venezuela-e-elogia-posicao-do-brasil-veja-video.ghtml"
$ws.Range("F94").Value = "lula"
$g94 = @"
idas militares ou sanções.
A nota oficial do Brasil, assinada pelo presidente Luiz Inácio Lula da Silva (PT) mencionou que os ataques “ultrapassam uma linha inaceitável” e configuram “
"@
$ws.Range("G94").Value = $g94
# The multi-line text above makes Excel mark the row with an explicit
# custom height; auto-fit it back so the row stays at the default height
# (matching how the rest of the sheet's multi-line cells are stored).
$ws.Rows(94).AutoFit()

# Row 95: new news item appended by the bot run.
$ws.Range("A95").Value = "05/01/2026 14:21:05"
$ws.Range("B95").Value = "05/01 13:59"
$ws.Range("C95").Value = "Metrópoles"
$ws.Range("D95").Value = "Vereador de Recife pede impeachment do prefeito João Campos"
$ws.Range("E95").Value = "https://www.metropoles.com/colunas/paulo-cappelli/vereador-de-recife-pede-impeachment-do-prefeito-joao-campos"
$ws.Range("F95").Value = "câmara"
$ws.Range("G95").Value = "r do Recife questiona nomeação em concurso da Procuradoria e pede abertura de processo na Câmara Municipal"
